$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph. We build the paragraph via InsertXML so we get
#    the exact run layout used elsewhere in the document: a leading empty
#    run, a bold "Meta description" run, and a plain run with the rest of
#    the sentence. A second, throw-away <w:p> is appended only to force a
#    paragraph break for the inserted content; it is deleted right after.
# ---------------------------------------------------------------------------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$metaXml = '<w:p ' + $wNs + '><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our Break Bones slot review to play for free &amp; discover how special features like Wild multipliers &amp; the Bonus Round can bring big wins.</w:t></w:r></w:p><w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Normal"/></w:pPr></w:p>'

$firstPara = $d.Paragraphs(1)
$insertionPoint = $d.Range($firstPara.Range.End, $firstPara.Range.End)
$insertionPoint.InsertXML($metaXml)

# Remove the placeholder empty paragraph that InsertXML needed in order to
# terminate the new paragraph (it is now Paragraphs(3), right before the
# original second paragraph).
$placeholder = $d.Paragraphs(3)
$d.Range($placeholder.Range.Start, $placeholder.Range.End).Delete()

# ---------------------------------------------------------------------------
# 2. Remove the paragraph near the end of the document that duplicated the
#    title ("Play Break Bones Free: Special Features & Big Wins - Review"
#    in bold) - that content now lives in the new meta-description
#    paragraph instead. It is always the second-to-last paragraph, right
#    before the italic call-to-action paragraph.
# ---------------------------------------------------------------------------
$titleOld = "Play Break Bones Free: Special Features & Big Wins - Review"
$count = $d.Paragraphs.Count
$titlePara = $d.Paragraphs($count - 1)
$titleText = $titlePara.Range.Text.TrimEnd([char]13, [char]7)
if ($titleText -eq $titleOld) {
    $d.Range($titlePara.Range.Start, $titlePara.Range.End).Delete()
} else {
    Write-Host "WARNING: expected duplicate title paragraph not found; got '$titleText'"
}

# ---------------------------------------------------------------------------
# 3. Swap out the italic call-to-action sentence for the new art-direction
#    prompt text. We set Range.Text directly (instead of Find/Replace's
#    replacement argument) so that straight apostrophes in the new text are
#    not silently "smart-quoted" into curly ones, and so we don't have to
#    worry about the same old sentence also appearing as a non-italic
#    fragment inside the freshly-inserted meta-description paragraph.
# ---------------------------------------------------------------------------
$oldItalic = "Read our Break Bones slot review to play for free & discover how special features like Wild multipliers & the Bonus Round can bring big wins."
$newItalic = "Please create a cartoon-style feature image for Break Bones that showcases a happy Maya warrior with glasses. The image should be eye-catching and draw in potential players, highlighting the game's theme and exciting features. Use bold, bright colors to make the image pop, and consider incorporating elements of the game, such as the Wilds and Scatter symbols, into the design. Make sure the Maya warrior is front and center, looking happy and excited to play the game. Overall, the image should convey a sense of fun and adventure, inviting players to join in on the action and give Break Bones a spin."

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastParaText = $lastPara.Range.Text.TrimEnd([char]13, [char]7)
if ($lastParaText -eq $oldItalic) {
    # Exclude the trailing paragraph-mark character from the replaced range.
    $contentRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
    $contentRange.Text = $newItalic
} else {
    Write-Host "WARNING: expected italic CTA paragraph not found; got '$lastParaText'"
}

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
